# Apply the "Added two Users and added visibility functionality" edit.
#
# Summary of per-sheet changes (see commit diff):
#   wheat   (sheet "wheat")   : drop the example data row and the
#                               State columns -> header-only A1:D1
#                               (SourceRailHead, DestinationRailHead,
#                               Commodity, Values)
#   rra     (sheet "rra")     : drop the State columns, replace the
#                               example row with BRCY/HSRA/RRA/1
#                               -> A1:D2
#   frk_rra (sheet "frk_rra") : drop the example data row -> header-only
#                               A1:F1
#   frk     (sheet "frk")     : add a new data row (AWB/Maharashtra/
#                               FCSJ/Jammu & Kashmir/FRK/1) -> A1:F2
#   frkcgr  (sheet "frkcgr")  : add two new data rows (JJKR/... and
#                               HKG/...) -> A1:F3

$wb = $excel.ActiveWorkbook

# --- wheat: remove sample row, collapse State columns out ---------------
$ws = $wb.Worksheets.Item("wheat")
$ws.Range("A2:F2").EntireRow.Delete()
$ws.Range("E1:F1").EntireColumn.Delete()
$ws.Range("A1").Value = "SourceRailHead"
$ws.Range("B1").Value = "DestinationRailHead"
$ws.Range("C1").Value = "Commodity"
$ws.Range("D1").Value = "Values"

# --- rra: collapse State columns out, replace sample row ----------------
$ws = $wb.Worksheets.Item("rra")
$ws.Range("E1:F2").EntireColumn.Delete()
$ws.Range("A1").Value = "SourceRailHead"
$ws.Range("B1").Value = "DestinationRailHead"
$ws.Range("C1").Value = "Commodity"
$ws.Range("D1").Value = "Values"
$ws.Range("A2").Value = "BRCY"
$ws.Range("B2").Value = "HSRA"
$ws.Range("C2").Value = "RRA"
$ws.Range("D2").Value = 1

# --- frk_rra: remove sample row -----------------------------------------
$ws = $wb.Worksheets.Item("frk_rra")
$ws.Range("A2:F2").EntireRow.Delete()

# --- frk: add new user row ----------------------------------------------
$ws = $wb.Worksheets.Item("frk")
$ws.Range("A2").Value = "AWB"
$ws.Range("B2").Value = "Maharashtra"
$ws.Range("C2").Value = "FCSJ"
$ws.Range("D2").Value = "Jammu & Kashmir"
$ws.Range("E2").Value = "FRK"
$ws.Range("F2").Value = 1

# --- frkcgr: add two new user rows --------------------------------------
$ws = $wb.Worksheets.Item("frkcgr")
$ws.Range("A2").Value = "JJKR"
$ws.Range("B2").Value = "Odisha"
$ws.Range("C2").Value = "BGTA"
$ws.Range("D2").Value = "MP"
$ws.Range("E2").Value = "FRK+CGR"
$ws.Range("F2").Value = 1
$ws.Range("A3").Value = "HKG"
$ws.Range("B3").Value = "Odisha"
$ws.Range("C3").Value = "BGTA"
$ws.Range("D3").Value = "MP"
$ws.Range("E3").Value = "FRK+CGR"
$ws.Range("F3").Value = 1
